# Refresh cryptos list price (D) and 1h volume-change (E) figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '58.988.13'
$ws.Range("E2").Value = '  -0.34%  '
$ws.Range("D3").Value = '2.515.90'
$ws.Range("E3").Value = '  +0.62%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = "'532.39"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.42%  '
$ws.Range("D6").Value = "'139.28"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.08%  '
$ws.Range("D7").Value = "'0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.15%  '
$ws.Range("D8").Value = "'0.563"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.51%  '
$ws.Range("D9").Value = '2.519.82'
$ws.Range("E10").Value = '  -0.05%  '
$ws.Range("E11").Value = '  +0.84%  '
$ws.Range("D12").Value = "'5.44"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -2.33%  '
$ws.Range("E13").Value = '  +0.40%  '
$ws.Range("D14").Value = '2.962.43'
$ws.Range("E14").Value = '  +0.64%  '
$ws.Range("D15").Value = "'23.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.56%  '
$ws.Range("D16").Value = '58.943.79'
$ws.Range("E16").Value = '  -0.22%  '
$ws.Range("E17").Value = '  +0.11%  '
$ws.Range("D18").Value = '2.515.51'
$ws.Range("E18").Value = '  -0.19%  '
$ws.Range("E19").Value = '  -1.43%  '
$ws.Range("D20").Value = "'4.26"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.24%  '
$ws.Range("D21").Value = "'322.77"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -0.20%  '
$ws.Range("E22").Value = '  +0.12%  '
$ws.Range("E23").Value = '  +1.22%  '
$ws.Range("D24").Value = "'62.14"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.35%  '
$ws.Range("D25").Value = "'0.426"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -3.03%  '
$ws.Range("E26").Value = '  +1.79%  '
$ws.Range("E27").Value = '  +0.52%  '
$ws.Range("D28").Value = "'7.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.23%  '
$ws.Range("D29").Value = "'6.73"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.70%  '
$ws.Range("D30").Value = '0.0₃0771'
$ws.Range("E30").Value = '  -0.45%  '
$ws.Range("D31").Value = "'1.78"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.20%  '
$ws.Range("D32").Value = "'163.38"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.79%  '
$ws.Range("E33").Value = '  +0.20%  '
$ws.Range("D34").Value = "'1.11"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -4.90%  '
$ws.Range("D35").Value = "'1.43"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.96%  '
$ws.Range("D36").Value = "'18.44"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.04%  '
$ws.Range("D37").Value = "'4.23"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.52%  '
$ws.Range("E38").Value = '  -1.99%  '
$ws.Range("D39").Value = "'36.88"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.31%  '
$ws.Range("E40").Value = '  -1.26%  '
$ws.Range("D41").Value = "'0.802"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.81%  '
$ws.Range("D42").Value = "'5.19"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -7.77%  '
$ws.Range("D43").Value = "'278.74"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -5.65%  '
$ws.Range("D44").Value = "'0.998"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.24%  '
$ws.Range("D45").Value = "'10.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.72%  '
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("E47").Value = '  +0.53%  '
$ws.Range("D48").Value = "'122.32"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.29%  '
$ws.Range("D49").Value = "'18.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.59%  '
$ws.Range("E50").Value = '  -0.57%  '
$ws.Range("E51").Value = '  -2.04%  '
